$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Long text blocks (new / moved content) ---

$objText = '1. Entendimento da relação entre a termodinâmica de soluções e os diagramas de fases.2. Domínio da leitura de diagramas unários, binários e ternários (configuração do sistema para um estado termodinâmico, leitura de composições de fases e cálculo de suas quantidades).3. Relacionamento entre microestruturas e diagramas de fases.4. Entendimento da seqüência de eventos que ocorrem no processo de solidificação em equilíbrio e fora de equilíbrio.'

$shortSyllabus = 'A. Introdução; teoria básica de equilíbrio de fases;B. Sistemas unários;C. Sistemas binários;D. Sistemas ternários;E. Cálculo termodinâmico de diagramas de fases.'

$fullSyllabus = '1. Introdução; revisão da termodinâmica de soluções; teoria básica de equilíbrio de fases; curvas de energia livre versus composição; regra das fases; 2. Sistemas unários, equilíbrios bi-, mono- e invariantes; 3. Sistemas binários isomorfos; a regra da alavanca; solidificação em equilíbrio e fora de equilíbrio; mínimos e máximos; 4. Sistemas eutéticos binários; solidificação e microetruturas de ligas hipoeutéticas, eutéticas e hipereutéticas; solidificação unidirecional com eutéticos; casos limites de eutéticos; 5. Sistemas eutetóides binários; solidificação e microetruturas de ligas hipoeutetóides, eutetói-des e hipereutetóides; o sistema Fe-C; 6. Sistemas monotéticos; sistemas monotetóides; sistemas metatéticos; transformações congruentes; 7. Sistemas peritéticos binários; resfriamento em equilíbrio e fora do equilíbrio de ligas peritéticas; sistemas peritetóides binários; sistemas sintéticos binários; 8. Sistemas ternários isomorfos; o triângulo de Gibbs; seções isotérmicas; projeções liquidus; seções verticais; máximos e mínimos; resfriamento em equilíbrio; 9. Equilíbrio ternário de três fases; regra da alavanca em campos trifásicos; resfriamento em equilíbrio; 10. Equilíbrio ternário de quatro fases: equilíbrio de classe I; equilíbrio de classe II e equilíbrio de classe III; 11. Transformações congruentes em sistemas ternários; sistemas ternários complexos; 12. Cálculo termodinâmico de diagramas de fases.'

$bibliography = '01. Gordon, P. Principles of Phase Diagrams in Materials Systems, McGraw-Hill, 1968.02. Rhines, F. N. Phase Diagrams in Metallurgy: Their Development and Applications, McGraw-Hill, 1956.03. Prince, A. Alloy Phase Equilibria, Elsevier, 1966.04. Massalski, T. B. Binary Alloys Phase Diagrams, ASM, Metals Park, Ohio, 1990.05. Alloy Phase Diagrams, ASM Handbook, Volume 3, ASM, Metals Park, Ohio, 1992.06. Hansen, M. Constitution of Binary Alloys, McGraw-Hill, 1958.07. Elliot, R. P. Constitution of Binary Alloys: First Supplement, McGraw-Hill, 1965.08. Shunk, F. A. Constitution of Binary Alloys: Second Supplement, McGraw-Hill, 1969.09. Levin, E. M. Phase Diagram for Ceramists, The American Ceramic Society, 1964.10. Rudman, P. S. Phase Stability in Metals and Alloys, McGraw-Hill, 1967.11. Kaufman, L. Computer Calculation of Phase Diagrams with Special Reference to Refractory Metals, Academic Press.12. Hack, K. The SGTE Casebook - Thermodynamics at Work. The Institut of Metals, London,6.13. Hillert, M. Phase Equilibria, Phase Diagrams and Phase Transformations. Cambridge University Press, Cambridge, 1998.14. Thermocalc version M manuals: User Guide and Examples, ThermoCalc AB, Stockholm, 1997."'

# --- Make room for the extra row: a new blank row is inserted at the current
#     row 13 ("Programa resumido:" / "Semestral"), pushing it (and everything
#     below) down by one and keeping its own height/style intact. The new
#     blank row 13 becomes the plain professor-name value row (no label),
#     and the old row (now row 14) is overwritten with the "Programa
#     resumido:" label + short-syllabus text. ---
$ws.Rows.Item(13).Insert()

# Row 10 ("Objetivos:") gets the long objectives text instead of the professor name.
$ws.Cells.Item(10, 2).Value = $objText
$ws.Cells.Item(10, 3).Value = $objText

# New row 13: just the professor-name value, no label in column A.
$ws.Cells.Item(13, 2).Value = '5009972 - Gilberto Carvalho Coelho'
$ws.Cells.Item(13, 3).Value = '5009972 - Gilberto Carvalho Coelho'

# Row 14 ("Programa resumido:", previously holding "Semestral") gets the short
# syllabus text.
$ws.Cells.Item(14, 2).Value = $shortSyllabus
$ws.Cells.Item(14, 3).Value = $shortSyllabus

# Row 16 ("Programa:") gets the full syllabus text instead of the activation date.
$ws.Cells.Item(16, 2).Value = $fullSyllabus
$ws.Cells.Item(16, 3).Value = $fullSyllabus

# Row 19 ("Método:") gets the "course will be taught..." method text instead of the
# professor name.
$ws.Cells.Item(19, 2).Value = 'O curso será ministrado na forma de aulas expositivas e aulas práticas em laboratório envolvendo preparação de amostras e caracterização microestrutural. Os resultados das aulas práticas serão apresentados oralmente e sujeitos a avaliação (T).'
$ws.Cells.Item(19, 3).Value = 'O curso será ministrado na forma de aulas expositivas e aulas práticas em laboratório envolvendo preparação de amostras e caracterização microestrutural. Os resultados das aulas práticas serão apresentados oralmente e sujeitos a avaliação (T).'

# Row 20 ("Critério:") gets the written-evaluations criteria text.
$ws.Cells.Item(20, 2).Value = 'Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF) juntamente com a avaliação do trabalho prático (T). O critério para a nota final é:NF=((P1*0,8)+(T*0,2)+P2*1)/2'
$ws.Cells.Item(20, 3).Value = 'Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF) juntamente com a avaliação do trabalho prático (T). O critério para a nota final é:NF=((P1*0,8)+(T*0,2)+P2*1)/2'

# Row 21 ("Norma de recuperação:") gets the recovery-grade norm text.
$ws.Cells.Item(21, 2).Value = 'Para os alunos que obtiverem 3,0≤NF<5,0, será aplicada uma avaliação de recuperação (R) que levará ao cálculo da média final (MF) com o seguinte critério:MF=(NF+R)/2'
$ws.Cells.Item(21, 3).Value = 'Para os alunos que obtiverem 3,0≤NF<5,0, será aplicada uma avaliação de recuperação (R) que levará ao cálculo da média final (MF) com o seguinte critério:MF=(NF+R)/2'

# Row 22 ("Bibliografia:") gets the bibliography text.
$ws.Cells.Item(22, 2).Value = $bibliography
$ws.Cells.Item(22, 3).Value = $bibliography

# --- Column layout fix: column B should be its own 60.71-wide range, matching
#     column C (previously column A's <col> definition incorrectly spanned
#     columns 1-2, so column B inherited column A's narrower width). ---
$ws.Columns.Item(2).ColumnWidth = 59.83
